$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Q0) updated values
$ws.Range("B3").Value = 0.2276211974233028
$ws.Range("C3").Value = 0.7533946750976762
$ws.Range("D3").Value = 1.127004305186878
$ws.Range("E3").Value = 1.061604589848253
$ws.Range("F3").Value = 1.041226601882469
$ws.Range("G3").Value = 121

# Row 4 (Q1) updated values
$ws.Range("B4").Value = 0.1031743366226042
$ws.Range("C4").Value = 0.5300950349840287
$ws.Range("D4").Value = 0.3665963335635029
$ws.Range("E4").Value = 0.6054719923856948
$ws.Range("F4").Value = 0.6017378854235643
$ws.Range("G4").Value = 59
